# Fruta / hortaliza, semanal
# Insert two new weekly rows (Feria Lagunitas de Puerto Montt - Durazno - Carson)
# above the previous row 208, shifting the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 208-209 (existing rows 208:241 shift down to 210:243).
$ws.Range("A208:A209").EntireRow.Insert()

# --- New row 208 ---
$ws.Range("A208").Value = 4
$ws.Range("B208").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C208").Value = "Los Lagos"
$ws.Range("D208").Value = 44943
$ws.Range("E208").Value = 10
$ws.Range("F208").Value = "Fruta"
$ws.Range("G208").Value = 100103
$ws.Range("H208").Value = "Frutos de hueso (carozo)"
$ws.Range("I208").Value = 100103004
$ws.Range("J208").Value = "Durazno"
$ws.Range("K208").Value = "Carson"
$ws.Range("L208").Value = "Especial"
$ws.Range("M208").Value = 200
$ws.Range("N208").Value = 23000
$ws.Range("O208").Value = 23000
$ws.Range("P208").Value = 23000
$ws.Range("Q208").Value = "$/caja 14 kilos empedrada"
$ws.Range("R208").Value = "Región de O'Higgins"
$ws.Range("S208").Value = 1643
$ws.Range("T208").Value = 14

# --- New row 209 ---
$ws.Range("A209").Value = 4
$ws.Range("B209").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C209").Value = "Los Lagos"
$ws.Range("D209").Value = 44943
$ws.Range("E209").Value = 10
$ws.Range("F209").Value = "Fruta"
$ws.Range("G209").Value = 100103
$ws.Range("H209").Value = "Frutos de hueso (carozo)"
$ws.Range("I209").Value = 100103004
$ws.Range("J209").Value = "Durazno"
$ws.Range("K209").Value = "Carson"
$ws.Range("L209").Value = "Primera"
$ws.Range("M209").Value = 400
$ws.Range("N209").Value = 19000
$ws.Range("O209").Value = 20000
$ws.Range("P209").Value = 19500
$ws.Range("Q209").Value = "$/caja 14 kilos empedrada"
$ws.Range("R209").Value = "Región de O'Higgins"
$ws.Range("S209").Value = 1393
$ws.Range("T209").Value = 14
